$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM values (new script run) for FAPs/MuSCs -> Wnt5b/Fzd8 edges.
# Keyed by worksheet column number -> new value for each data row.
$rowUpdates = @{
  2 = @{ 9=0.8757151715610434; 10=0.8757151715610434; 13=2.544438666666667; 14=7.633316000000001; 15=0.201325300207035; 16=0.201325300207035; 17=0.6024832133942223; 18=5.422348920548001; 19=0.1763036198103822; 20=0.1763036198103822 }
  3 = @{ 9=0.8757151715610434; 10=0.8757151715610434; 15=0.6969390273602759; 16=0.696939027360276; 19=0.6103200799123907; 20=0.6103200799123908 }
  4 = @{ 4="MuSCs"; 9=0.8757151715610434; 10=0.8757151715610434; 11=3; 12=1; 13=1.273916333333333; 14=3.821749; 15=0.1007969229547075; 16=0.1007969229547075; 17=0.3016434297107778; 18=2.714790867397; 19=0.08826939467810692; 20=0.08826939467810693 }
  5 = @{ 4="Resolving-Mac"; 9=0.8757151715610434; 10=0.8757151715610434; 11=1; 12=0.3333333333333333; 13=0.01186433333333333; 14=0.035593; 15=0.0009387494779816524; 16=0.0009387494779816526; 17=0.002809288258777778; 18=0.025283594329; 19=0.0008220771601635426; 20=0.0008220771601635429 }
  6 = @{ 1="MuSCs"; 4="ECs"; 7=0.03360533333333333; 8=0.100816; 9=0.1242848284389566; 10=0.1242848284389566; 11=3; 12=1; 13=2.544438666666667; 14=7.633316000000001; 15=0.201325300207035; 16=0.201325300207035; 17=0.08550670953955557; 18=0.7695603858560001; 19=0.02502168039665278; 20=0.02502168039665278 }
  7 = @{ 4="FAPs"; 7=0.03360533333333333; 8=0.100816; 9=0.1242848284389566; 10=0.1242848284389566; 13=8.808225333333333; 14=26.424676; 15=0.6969390273602759; 16=0.696939027360276; 17=0.2960033484017778; 18=2.664030135616; 19=0.08661894744788518; 20=0.0866189474478852 }
  8 = @{ 4="MuSCs"; 7=0.03360533333333333; 8=0.100816; 9=0.1242848284389566; 10=0.1242848284389566; 13=1.273916333333333; 14=3.821749; 15=0.1007969229547075; 16=0.1007969229547075; 17=0.04281038302044445; 18=0.385293447184; 19=0.01252752827660054; 20=0.01252752827660055 }
  9 = @{ 4="Resolving-Mac"; 7=0.03360533333333333; 8=0.100816; 9=0.1242848284389566; 10=0.1242848284389566; 13=0.01186433333333333; 14=0.035593; 15=0.0009387494779816524; 16=0.0009387494779816526; 17=0.0003987048764444444; 18=0.003588343888; 19=0.0001166723178181097; 20=0.0001166723178181098 }
}

foreach ($rowNum in $rowUpdates.Keys) {
  $colUpdates = $rowUpdates[$rowNum]
  foreach ($colNum in $colUpdates.Keys) {
    $ws.Cells.Item($rowNum, $colNum).Value = $colUpdates[$colNum]
  }
}

# The refreshed TPM run collapsed the duplicate "Inflammatory-Mac" target-cluster
# rows (old rows 10 and 11) into the recomputed rows above, so remove them and
# let the sheet shrink from A1:T11 to A1:T9.
$ws.Rows(10).Delete()
$ws.Rows(10).Delete()
